$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '88.467.62'
$ws.Range("E2").Value = '  +9.69%  '
$ws.Range("D3").Value = '3.330.35'
$ws.Range("E3").Value = '  +6.69%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.75'
$ws.Range("E5").Value = '  +6.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '650.42'
$ws.Range("E6").Value = '  +4.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.400'
$ws.Range("E7").Value = '  +44.38%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.602'
$ws.Range("E9").Value = '  +5.02%  '
$ws.Range("D10").Value = '3.327.67'
$ws.Range("E10").Value = '  +6.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.587'
$ws.Range("E11").Value = '  +3.15%  '
$ws.Range("E12").Value = '  +17.63%  '
$ws.Range("E13").Value = '  +2.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.61'
$ws.Range("E14").Value = '  +15.29%  '
$ws.Range("D15").Value = '3.936.17'
$ws.Range("E15").Value = '  +6.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.52'
$ws.Range("E16").Value = '  +5.76%  '
$ws.Range("D17").Value = '88.307.94'
$ws.Range("E17").Value = '  +9.11%  '
$ws.Range("D18").Value = '3.326.14'
$ws.Range("E18").Value = '  +6.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.65'
$ws.Range("E19").Value = '  +6.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.12'
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.61'
$ws.Range("E21").Value = '  +8.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '455.68'
$ws.Range("E22").Value = '  +7.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.45'
$ws.Range("E23").Value = '  +8.97%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.46'
$ws.Range("E24").Value = '  +4.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.60'
$ws.Range("E25").Value = '  +10.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.97'
$ws.Range("E26").Value = '  +21.53%  '
$ws.Range("D27").Value = '3.494.52'
$ws.Range("E27").Value = '  +5.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '78.15'
$ws.Range("E28").Value = '  +4.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.212'
$ws.Range("E29").Value = '  +46.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0000134'
$ws.Range("E30").Value = '  +13.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.40'
$ws.Range("E32").Value = '  +6.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '598.90'
$ws.Range("E33").Value = '  +9.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.60'
$ws.Range("E34").Value = '  +10.13%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  +7.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.24'
$ws.Range("E37").Value = '  +24.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.142'
$ws.Range("E38").Value = '  -4.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.15'
$ws.Range("E39").Value = '  +3.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.17'
$ws.Range("E40").Value = '  +10.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.421'
$ws.Range("E41").Value = '  +5.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.85'
$ws.Range("E42").Value = '  +5.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.17'
$ws.Range("E44").Value = '  +6.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '158.50'
$ws.Range("E45").Value = '  -0.98%  '
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.45'
$ws.Range("E47").Value = '  +11.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '187.67'
$ws.Range("E48").Value = '  +1.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.96'
$ws.Range("E49").Value = '  +6.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.46'
$ws.Range("E50").Value = '  +7.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.662'
$ws.Range("E51").Value = '  +7.72%  '
